# Reduce the "custom labels" sort fixture down to a smaller 4-row sample:
# - Binary sheet: keep header + first 3 data rows, delete the rest.
# - Ternary / Quaternary sheets: keep header + first 3 data rows, blank out
#   the remaining (already-present) rows instead of deleting them.
# - Update the active sheet / selections to match the new, smaller ranges.

$wb = $excel.ActiveWorkbook

$wsBinary = $wb.Worksheets.Item("Binary")
$wsTernary = $wb.Worksheets.Item("Ternary")
$wsQuaternary = $wb.Worksheets.Item("Quaternary")

# --- Binary: physically remove rows 5:10, leaving a 4-row (1 header + 3 data) table ---
[void]$wsBinary.Rows("5:10").Delete()
[void]$wsBinary.Range("A5:C11").Select()

# --- Quaternary: clear the now-unused data rows 5:20 ---
# (Note: Range.Select() implicitly activates its sheet, so this sheet's work
# happens before Ternary's, since Ternary must end up the active tab.)
[void]$wsQuaternary.Range("A5:D20").ClearContents()
[void]$wsQuaternary.Range("A5:D14").Select()
[void]$wsQuaternary.Range("D14").Activate()

# --- Ternary: clear the now-unused data rows 5:20, then make it the active sheet ---
[void]$wsTernary.Range("A5:C20").ClearContents()
[void]$wsTernary.Activate()
[void]$wsTernary.Range("A5:D25").Select()
